$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new columns of user detail data: first name, last name, pin code
$ws.Range("C1").Value = "first name"
$ws.Range("D1").Value = "last name"
$ws.Range("E1").Value = "pin code"

for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 3).Value = "fhjag"
    $ws.Cells.Item($r, 4).Value = "hfajh"
    $ws.Cells.Item($r, 5).Value = "gfhaf"
}

# Page setup (portrait, letter/A4-class paper) as left by the save
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection lands on the last filled cell
$ws.Range("E5").Select()
